# Apply the RPA datasets push 2023-10-24 update to Sheet1.
# The underlying data table (rows 2-27) is replaced with a re-sorted,
# 2-row-shorter table (rows 2-25): "대신밸런스제15호스팩" and "한국제12호스팩"
# rows are removed, and rows are re-ordered by subscription date within
# each underwriter group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (청약일), F (납입일), G (상장일) hold date-like text ("YYYY-MM-DD").
# Force them to Text format before writing so Excel does not silently
# reinterpret the strings as date serial numbers.
$ws.Range("B2:B25").NumberFormat = "@"
$ws.Range("F2:F25").NumberFormat = "@"
$ws.Range("G2:G25").NumberFormat = "@"

$data = @(
    @("CS","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",42120,16200000,26000,0,10),
    @("KB","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",42120,16200000,26000,0,10),
    @("KB","2023-09-19","한싹","KB","KB","2023-09-22","2023-10-04",18750,1500000,12500,0,100),
    @("NH","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",42120,16200000,26000,0,10),
    @("대신","2023-08-23","대신밸런스제16호스팩","대신","대신","2023-08-28","2023-09-04",13000,6500000,2000,0,100),
    @("미래","2023-09-18","밀리의서재","미래","미래","2023-09-21","2023-09-27",34500,1500000,23000,0,100),
    @("미래","2023-10-10","신성에스티","미래","미래","2023-10-13","2023-10-19",52000,2000000,26000,0,100),
    @("미래","2023-10-05","퓨릿","미래","미래","2023-10-11","2023-10-18",44265.9,4137000,10700,0,100),
    @("미래","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",126360,16200000,26000,0,30),
    @("삼성","2023-09-19","레뷰코퍼레이션","삼성","삼성","2023-09-22","2023-10-06",33600,2240000,15000,0,100),
    @("상상인","2023-09-04","상상인제4호스팩","상상인","상상인","2023-09-07","2023-09-14",9000,4500000,2000,0,100),
    @("신영","2023-09-14","인스웨이브시스템즈","신영","신영","2023-09-19","2023-09-25",26400,1100000,24000,0,100),
    @("신영","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",12636,16200000,26000,0,3),
    @("신한","2023-09-19","신한제11호스팩","신한","신한","2023-09-22","2023-10-04",36000,18000000,2000,0,100),
    @("유비에스","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",4212,16200000,26000,0,1),
    @("유안타","2023-09-18","아이엠티","유안타","유안타, 유진","2023-09-21","2023-10-10",15484,1580000,14000,0,70),
    @("유안타","2023-08-22","유안타제11호스팩","유안타","유안타","2023-08-25","2023-09-01",10000,5000000,2000,0,100),
    @("유진","2023-09-18","아이엠티","유안타","유안타, 유진","2023-09-21","2023-10-10",6636,1580000,14000,0,30),
    @("키움","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",12636,16200000,26000,0,3),
    @("하나","2023-10-10","에스엘에스바이오","하나","하나","2023-10-13","2023-10-20",5390,770000,7000,0,100),
    @("하나","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",12636,16200000,26000,0,3),
    @("한국","2023-09-21","두산로보틱스","한국, 미래","한국, 미래, NH, KB, CS, 키움, 신영, 하나, 유비에스리미티드(영업소)","2023-09-26","2023-10-05",126360,16200000,26000,0,30),
    @("한화","2023-08-29","한화플러스제4호스팩","한화","한화","2023-09-01","2023-09-07",9500,4750000,2000,0,100),
    @("현대차","2023-09-25","에이치엠씨제6호스팩","현대차","현대차","2023-10-04","2023-10-13",8000,4000000,2000,0,100)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value  = $row[4]
    $ws.Cells.Item($r, 6).Value  = $row[5]
    $ws.Cells.Item($r, 7).Value  = $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $r = $r + 1
}

# The table shrank from 27 data rows (2-27) to 24 (2-25); drop the two
# now-unused trailing rows so the sheet dimension becomes A1:L25.
$ws.Range("26:27").Delete()

